# Updated symbol list on Thu Dec 22 16:08:31 UTC 2022 with GitHub Actions
# Applies refreshed coin prices/ranks pulled from coinranking.com to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. Values are written as plain text
# (NumberFormat "@") so numeric-looking strings (prices, the hour column, ...)
# keep their exact original formatting (e.g. trailing zeros) instead of being
# coerced into floating point numbers by Excel. Style is restored to "Normal"
# afterwards so only the cell content changes.
$updates = [ordered]@{
    "D2" = "241.55"
    "G2" = "16"
    "D3" = "21.77"
    "G3" = "16"
    "D4" = "5.372"
    "G4" = "16"
    "D5" = "0.05688"
    "G5" = "16"
    "D6" = "3.419"
    "G6" = "16"
    "D7" = "6.289"
    "G7" = "16"
    "D8" = "0.8061"
    "G8" = "16"
    "D9" = "0.8369"
    "G9" = "16"
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D10" = "0.1426"
    "E10" = "9WazirXWRX"
    "G10" = "16"
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D11" = "0.07270"
    "E11" = "10MandalaExchangeTokenMDX"
    "G11" = "16"
    "B12" = "LiechtensteinCryptoassetsExchange"
    "C12" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D12" = "0.03043"
    "E12" = "11LiechtensteinCryptoassetsExchangeLCX"
    "G12" = "16"
    "B13" = "BitrueCoin"
    "C13" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D13" = "0.03167"
    "E13" = "12BitrueCoinBTR"
    "G13" = "16"
    "B14" = "BitMartToken"
    "C14" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D14" = "0.09346"
    "E14" = "13BitMartTokenBMX"
    "G14" = "16"
    "B15" = "MCDex"
    "C15" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "D15" = "3.918"
    "E15" = "14MCDexMCB"
    "G15" = "16"
    "B16" = "BitForexToken"
    "C16" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D16" = "0.001584"
    "E16" = "15BitForexTokenBF"
    "G16" = "16"
    "B17" = "CoinExToken"
    "C17" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "D17" = "0.04821"
    "E17" = "16CoinExTokenCET"
    "G17" = "16"
    "B18" = "One"
    "C18" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "D18" = "0.0005858"
    "E18" = "17OneONE"
    "G18" = "16"
    "D19" = "0.006344"
    "G19" = "16"
    "D20" = "0.0009948"
    "G20" = "16"
    "D21" = "0.004048"
    "G21" = "16"
    "G22" = "16"
    "D23" = "3.719"
    "G23" = "16"
    "D24" = "2.170"
    "G24" = "16"
    "D25" = "0.3232"
    "G25" = "16"
    "G26" = "16"
    "G27" = "16"
    "G28" = "16"
    "G29" = "16"
    "G30" = "16"
    "G31" = "16"
    "G32" = "16"
    "G33" = "16"
    "G34" = "16"
    "G35" = "16"
    "G36" = "16"
    "G37" = "16"
    "G38" = "16"
    "G39" = "16"
    "G40" = "16"
    "D41" = "0.006739"
    "G41" = "16"
    "G42" = "16"
    "G43" = "16"
    "D44" = "0.006562"
    "G44" = "16"
    "D45" = "0.00005616"
    "G45" = "16"
    "G46" = "16"
    "D47" = "0.5808"
    "E47" = "46CoinbaseStockTokenCOINBestin24h"
    "G47" = "16"
    "D48" = "0.1421"
    "G48" = "16"
    "D49" = "0.00002102"
    "G49" = "16"
    "G50" = "16"
    "G51" = "16"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
